# The "settings" sheet (XLSForm settings tab) gains a "version" column so
# tests can control the form version: header "version" in C1 and the
# numeric value 1 in C2, next to the existing form_title/form_id columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

$ws.Range("C1").Value = "version"
$ws.Range("C2").Value = 1

# Reflect the author's final cursor position (cell below the new column).
[void]$ws.Range("C3").Select()
